$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume/name/link cells keep their original text formatting
# (avoid Excel auto-converting numeric-looking strings like "1.010" into numbers)
$ws.Range("B2:E51").NumberFormat = "@"

$updates = @(
    ,@('D2', '27.872.46')
    ,@('E2', '  +0.15%  ')
    ,@('D3', '1.887.23')
    ,@('E3', '  +0.03%  ')
    ,@('D4', '1.012')
    ,@('D5', '337.08')
    ,@('E5', '  +0.92%  ')
    ,@('D6', '1.010')
    ,@('E6', '  +0.34%  ')
    ,@('D7', '0.4685')
    ,@('E7', '  -0.94%  ')
    ,@('D8', '0.3967')
    ,@('E8', '  +0.98%  ')
    ,@('D9', '45.93')
    ,@('E9', '  -3.80%  ')
    ,@('D10', '0.08026')
    ,@('E10', '  -0.43%  ')
    ,@('D11', '1.015')
    ,@('E11', '  -1.15%  ')
    ,@('D12', '22.04')
    ,@('E12', '  -0.26%  ')
    ,@('D13', '1.890.61')
    ,@('E13', '  +0.78%  ')
    ,@('D14', '6.016')
    ,@('E14', '  +0.43%  ')
    ,@('D15', '7.299')
    ,@('E15', '  +2.15%  ')
    ,@('D16', '1.013')
    ,@('E16', '  +0.30%  ')
    ,@('D17', '89.54')
    ,@('E17', '  +2.53%  ')
    ,@('D18', '0.06738')
    ,@('E18', '  +0.25%  ')
    ,@('D19', '0.00001049')
    ,@('E19', '  -0.22%  ')
    ,@('D20', '17.41')
    ,@('E20', '  +0.45%  ')
    ,@('D21', '1.009')
    ,@('E21', '  +0.23%  ')
    ,@('D22', '27.870.28')
    ,@('E22', '  +0.05%  ')
    ,@('D23', '5.511')
    ,@('E23', '  -0.24%  ')
    ,@('D24', '11.03')
    ,@('E24', '  +0.35%  ')
    ,@('D25', '2.317')
    ,@('E25', '  -0.58%  ')
    ,@('D26', '2.108.95')
    ,@('E26', '  +0.36%  ')
    ,@('D27', '159.31')
    ,@('E27', '  +0.14%  ')
    ,@('D28', '19.87')
    ,@('E28', '  -1.34%  ')
    ,@('D29', '2.164')
    ,@('E29', '  +2.77%  ')
    ,@('D30', '5.519')
    ,@('E30', '  -0.99%  ')
    ,@('D31', '122.17')
    ,@('E31', '  +0.08%  ')
    ,@('D32', '0.9880')
    ,@('E32', '  +1.01%  ')
    ,@('D33', '0.09490')
    ,@('E33', '  -0.13%  ')
    ,@('D34', '3.645')
    ,@('E34', '  +0.59%  ')
    ,@('D35', '5.354')
    ,@('E35', '  -0.07%  ')
    ,@('D36', '1.361')
    ,@('E36', '  -6.19%  ')
    ,@('E37', '  -1.11%  ')
    ,@('D38', '0.02251')
    ,@('E38', '  -0.90%  ')
    ,@('D39', '8.373')
    ,@('E39', '  +3.67%  ')
    ,@('D40', '1.205')
    ,@('E40', '  -1.19%  ')
    ,@('D41', '1.009')
    ,@('E41', '  +0.21%  ')
    ,@('D42', '0.6010')
    ,@('E42', '  +0.05%  ')
    ,@('D43', '0.1899')
    ,@('E43', '  +0.04%  ')
    ,@('D44', '10.43')
    ,@('E44', '  +1.18%  ')
    ,@('B45', 'WEMIXTOKEN')
    ,@('C45', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix')
    ,@('D45', '1.249')
    ,@('E45', '  -0.76%  ')
    ,@('B46', 'Decentraland')
    ,@('C46', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana')
    ,@('D46', '0.5677')
    ,@('E46', '  -0.60%  ')
    ,@('E47', '  +1.22%  ')
    ,@('D48', '1.952')
    ,@('E48', '  +0.56%  ')
    ,@('D49', '0.06791')
    ,@('E49', '  -1.76%  ')
    ,@('D50', '112.86')
    ,@('E50', '  -0.38%  ')
    ,@('D51', '3.037')
    ,@('E51', '  -10.68%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
